# Update countries & provincias Spain
#
# Refreshes the COVID-19 "Pais" sheet with newer source numbers:
#   - Belgica's case count overtakes Marruecos/Bolivia, so it now sorts
#     ahead of them (rows 33-35).
#   - Honduras's case count overtakes Venezuela, so it now sorts ahead of
#     it (rows 55-56).
#   - Several other countries (Pakistan, Kazajistan, San Martin (Parte
#     Holandesa), San Cristobal y Nieves) get refreshed totals.
#   - The "last updated" footer timestamp moves from 04:09 to 05:26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados..." footer timestamp -----------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 05:26"

# --- Row 25 : Pakistan ------------------------------------------------------
$ws.Range("B25").Value = 317595
$ws.Range("C25").Value = 661
$ws.Range("D25").Value = 302708
$ws.Range("E25").Value = 8335
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 6552

# --- Rows 33-35: Belgica overtakes Marruecos and Bolivia -------------------
$ws.Range("A33").Value = "Belgica"
$ws.Range("B33").Value = 143596
$ws.Range("C33").Value = 5728
$ws.Range("D33").Value = 19981
$ws.Range("E33").Value = 113489
$ws.Range("G33").Value = 18
$ws.Range("H33").Value = 10126

$ws.Range("A34").Value = "Marruecos"
$ws.Range("B34").Value = 142953
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 120275
$ws.Range("E34").Value = 20192
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 2486

$ws.Range("A35").Value = "Bolivia"
$ws.Range("B35").Value = 137969
$ws.Range("C35").Value = 263
$ws.Range("D35").Value = 99812
$ws.Range("E35").Value = 29929
$ws.Range("G35").Value = 36
$ws.Range("H35").Value = 8228

# --- Row 41: Kazajistan ------------------------------------------------------
$ws.Range("B41").Value = 108561
$ws.Range("C41").Value = 107
$ws.Range("D41").Value = 103758
$ws.Range("E41").Value = 3057

# --- Rows 55-56: Honduras overtakes Venezuela -------------------------------
$ws.Range("A55").Value = "Honduras"
$ws.Range("B55").Value = 81672
$ws.Range("C55").Value = 512
$ws.Range("D55").Value = 31089
$ws.Range("E55").Value = 48106
$ws.Range("G55").Value = 11
$ws.Range("H55").Value = 2477

$ws.Range("A56").Value = "Venezuela"
$ws.Range("B56").Value = 81019
$ws.Range("D56").Value = 72196
$ws.Range("E56").Value = 8145
$ws.Range("H56").Value = 678

# --- Row 172: San Martin (Parte Holandesa) ----------------------------------
$ws.Range("B172").Value = 699
$ws.Range("C172").Value = 1
$ws.Range("D172").Value = 620
$ws.Range("E172").Value = 57

# --- Row 211: San Cristobal y Nieves ----------------------------------------
$ws.Range("D211").Value = 18
$ws.Range("E211").Value = 1
